$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.365.10"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.280.82"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'112.38"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").Value = "'265.47"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "'0.651"
$ws.Range("E7").Value = "  +4.09%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.612"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'46.85"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").Value = "'0.0937"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "'9.28"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "'15.31"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "2.617.60"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "'0.865"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "2.275.18"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "43.193.96"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "'0.0000109"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "'6.74"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "'72.19"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'2.43"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "'234.93"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'2.90"
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("D25").Value = "'9.37"
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").Value = "'11.43"
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("D28").Value = "'41.06"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").Value = "'3.34"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").Value = "'173.17"
$ws.Range("D32").Value = "'21.67"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "'0.0897"
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").Value = "'5.64"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "'4.69"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "'3.89"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").Value = "'2.59"
$ws.Range("E40").Value = "  +7.66%  "
$ws.Range("D41").Value = "'14.34"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("D42").Value = "'74.31"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").Value = "'0.237"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("D44").Value = "'6.09"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "'1.27"
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("D48").Value = "'8.56"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "'0.0997"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "'100.40"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").Value = "'0.607"
$ws.Range("E51").Value = "  +11.35%  "
